# Auto-generated Excel COM-interop script applying scheduled market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 186
$ws.Range("I38").Value = 186
$ws.Range("K38").Value = 558
$ws.Range("M38").Value = -186
$ws.Range("H58").Value = 504
$ws.Range("J58").Value = 1500
$ws.Range("L58").Value = 4500
$ws.Range("N58").Value = -4800
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H86").Value = 131500.75
$ws.Range("I86").Value = 1003
$ws.Range("J86").Value = 175000
$ws.Range("K86").Value = 1003
$ws.Range("L86").Value = 175000
$ws.Range("M86").Value = 120
$ws.Range("N86").Value = -177246
$ws.Range("H89").Value = 131500.75
$ws.Range("I89").Value = 1003
$ws.Range("J89").Value = 175000
$ws.Range("K89").Value = 5015
$ws.Range("L89").Value = 875000
$ws.Range("M89").Value = 601
$ws.Range("N89").Value = -886232
$ws.Range("H113").Value = 8106.143
$ws.Range("J113").Value = 8998.375
$ws.Range("L113").Value = 8998.375
$ws.Range("N113").Value = -15506.375
$ws.Range("H132").Value = 1000
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 3000
$ws.Range("N132").Value = -8060
$ws.Range("H134").Value = 780000
$ws.Range("J134").Value = 780000
$ws.Range("L134").Value = 780000
$ws.Range("N134").Value = -790140
$ws.Range("H135").Value = 891.5
$ws.Range("I135").Value = 869.8
$ws.Range("K135").Value = 7828.2
$ws.Range("M135").Value = -5293.2
$ws.Range("H137").Value = 4908
$ws.Range("I137").Value = 4837
$ws.Range("K137").Value = 14511
$ws.Range("M137").Value = -11961
$ws.Range("H138").Value = 2777.7778
$ws.Range("I138").Value = 2333.3333
$ws.Range("K138").Value = 6999.999899999999
$ws.Range("M138").Value = -1859.999899999999
$ws.Range("H141").Value = 2219.889
$ws.Range("I141").Value = 2247.625
$ws.Range("K141").Value = 6742.875
$ws.Range("M141").Value = -1562.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H32").Value = 2260.862
$ws.Range("I32").Value = 1687.6296
$ws.Range("K32").Value = 1687.6296
$ws.Range("M32").Value = -1400.6296
$ws.Range("H63").Value = 2068.25
$ws.Range("I63").Value = 1757.6666
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1757.6666
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1071.6666
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2068.25
$ws.Range("I66").Value = 1757.6666
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 8788.333000000001
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -5356.333000000001
$ws.Range("N66").Value = -21864
$ws.Range("H110").Value = 2434.5
$ws.Range("I110").Value = 2477.8572
$ws.Range("K110").Value = 2477.8572
$ws.Range("M110").Value = -432.8571999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 251
$ws.Range("I11").Value = 251
$ws.Range("K11").Value = 251
$ws.Range("M11").Value = -111
$ws.Range("H24").Value = 516
$ws.Range("I24").Value = 516
$ws.Range("K24").Value = 516
$ws.Range("M24").Value = -281
$ws.Range("H82").Value = 19150.334
$ws.Range("H85").Value = 19150.334
$ws.Range("H94").Value = 2437.3
$ws.Range("I94").Value = 2595.889
$ws.Range("K94").Value = 2595.889
$ws.Range("M94").Value = -2144.889
$ws.Range("H134").Value = 16666.334
$ws.Range("J134").Value = 14999
$ws.Range("L134").Value = 44997
$ws.Range("N134").Value = -50067
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2291.3462
$ws.Range("I31").Value = 1942.1538
$ws.Range("J31").Value = 2640.5386
$ws.Range("K31").Value = 1942.1538
$ws.Range("L31").Value = 2640.5386
$ws.Range("M31").Value = -1647.1538
$ws.Range("N31").Value = -3230.5386
$ws.Range("H34").Value = 2291.3462
$ws.Range("I34").Value = 1942.1538
$ws.Range("J34").Value = 2640.5386
$ws.Range("K34").Value = 1942.1538
$ws.Range("L34").Value = 2640.5386
$ws.Range("M34").Value = -1740.1538
$ws.Range("N34").Value = -3044.5386
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1418.4
$ws.Range("J5").Value = 1648.5
$ws.Range("L5").Value = 4945.5
$ws.Range("N5").Value = -5169.5
$ws.Range("H11").Value = 500100.5
$ws.Range("J11").Value = 151
$ws.Range("L11").Value = 453
$ws.Range("N11").Value = -733
$ws.Range("H16").Value = 396.33334
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H32").Value = 894.5
$ws.Range("J32").Value = 894.5
$ws.Range("L32").Value = 2683.5
$ws.Range("N32").Value = -3249.5
$ws.Range("H128").Value = 308284.16
$ws.Range("I128").Value = 308284.16
$ws.Range("K128").Value = 924852.48
$ws.Range("M128").Value = -919872.48
$ws.Range("H131").Value = 1432.875
$ws.Range("I131").Value = 1096.6
$ws.Range("J131").Value = 1993.3334
$ws.Range("K131").Value = 3289.8
$ws.Range("L131").Value = 5980.0002
$ws.Range("M131").Value = 1750.2
$ws.Range("N131").Value = -16060.0002
$ws.Range("H135").Value = 1418.4
$ws.Range("J135").Value = 1648.5
$ws.Range("L135").Value = 14836.5
$ws.Range("N135").Value = -19906.5
$ws.Range("H137").Value = 1547.8334
$ws.Range("J137").Value = 2500
$ws.Range("L137").Value = 7500
$ws.Range("N137").Value = -17700
$ws.Range("H139").Value = 1983.4445
$ws.Range("I139").Value = 1481.75
$ws.Range("K139").Value = 4445.25
$ws.Range("M139").Value = 694.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4417.143
$ws.Range("I97").Value = 2970
$ws.Range("J97").Value = 5502.5
$ws.Range("K97").Value = 2970
$ws.Range("L97").Value = 5502.5
$ws.Range("M97").Value = -2474
$ws.Range("N97").Value = -6494.5
$ws.Range("H102").Value = 2846.9375
$ws.Range("I102").Value = 2789.0715
$ws.Range("K102").Value = 2789.0715
$ws.Range("M102").Value = -1167.0715
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24784
$ws.Range("I7").Value = 22911.5
$ws.Range("J7").Value = 27780
$ws.Range("K7").Value = 22911.5
$ws.Range("L7").Value = 27780
$ws.Range("M7").Value = -22799.5
$ws.Range("N7").Value = -28004
$ws.Range("H22").Value = 773.2273
$ws.Range("I22").Value = 710.55
$ws.Range("K22").Value = 710.55
$ws.Range("M22").Value = -415.55
$ws.Range("H27").Value = 773.2273
$ws.Range("I27").Value = 710.55
$ws.Range("K27").Value = 710.55
$ws.Range("M27").Value = -603.55
$ws.Range("H82").Value = 1112.8
$ws.Range("I82").Value = 1112.8
$ws.Range("K82").Value = 1112.8
$ws.Range("M82").Value = -751.8
$ws.Range("H85").Value = 1112.8
$ws.Range("I85").Value = 1112.8
$ws.Range("K85").Value = 1112.8
$ws.Range("M85").Value = 135.2
$ws.Range("H126").Value = 24784
$ws.Range("I126").Value = 22911.5
$ws.Range("J126").Value = 27780
$ws.Range("K126").Value = 68734.5
$ws.Range("L126").Value = 83340
$ws.Range("M126").Value = -66264.5
$ws.Range("N126").Value = -88280
$ws.Range("H136").Value = 3688
$ws.Range("J136").Value = 5946.5
$ws.Range("L136").Value = 17839.5
$ws.Range("N136").Value = -22939.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2745.111
$ws.Range("I136").Value = 2745.111
$ws.Range("K136").Value = 8235.332999999999
$ws.Range("M136").Value = -5685.332999999999

Write-Output "Applied 200 value updates and 4 cell clears"
